# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet and both locale sheets.
# - The locale sheets (zh-cn, de-de) get their "Latest Target File" (F) and
#   "Latest Handback File" (G) columns populated with hyperlinks, and their
#   "Latest Handback DateTime" (H) stamped with the actual handback time.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# --- Overview sheet: Status shown in columns B (zh-cn) and C (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Range("F2").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", "", "", "a.md") | Out-Null

$wsZh.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e089fe9d48f3e213d316fda9a5919c65c531736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

$wsZh.Range("F3").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", "", "", "a.md") | Out-Null

$wsZh.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e089fe9d48f3e213d316fda9a5919c65c531736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

$wsZh.Range("H2").Value = "2016-03-21 02:26:46"
$wsZh.Range("H3").Value = "2016-03-21 02:26:46"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Range("F2").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", "", "", "a.md") | Out-Null

$wsDe.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f43a641668f0430fd1b0a8146d4641ce6feb246/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$wsDe.Range("F3").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/fe09c337f8da52cbaad6a4687bf0dc7d45316f0b/e2e/a.md", "", "", "a.md") | Out-Null

$wsDe.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f43a641668f0430fd1b0a8146d4641ce6feb246/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$wsDe.Range("H2").Value = "2016-03-21 02:26:53"
$wsDe.Range("H3").Value = "2016-03-21 02:26:53"

Write-Host "Handback report generated."
